# Implemented Unreachable Activities test.
# Replace the old "Undocumented Comment Out activity" checklist row (row 12)
# with the new "Unreachable activities" check, and drop the now-superfluous
# trailing blank row (row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 12 with the new check's title and xaml filename.
# (Explanation/Suggestion text in columns E/F already describe unreachable
# activities and stay as-is.)
$ws.Cells.Item(12, 2).Value2 = "Unreachable activities"
$ws.Cells.Item(12, 3).Value2 = "Checks\UnreachableActivities.xaml"

# Remove the trailing empty row.
$ws.Rows.Item(13).Delete()

# Reflect the new selection left behind in the source workbook.
$ws.Range("C12").Select()
